# Apply the EDT_dict.xlsx update:
#  - Update the German (B6) and English (C6) ITEM_INSTRUCTION strings:
#      "Summton" -> "Piepton" (German)
#      "buzzer"  -> "beep"    (English, appears twice)
#  - Resize column B (narrower) and give column C an explicit width
#  - Update the saved selection to C6 (and drop the scrolled topLeftCell)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the ITEM_INSTRUCTION row text (row 6: DE in B, EN in C) ---
$deOld = $ws.Range("B6").Value2
$enOld = $ws.Range("C6").Value2

$deNew = $deOld.Replace("Summton", "Piepton")
$enNew = $enOld.Replace("buzzer", "beep")

$ws.Range("B6").Value = $deNew
$ws.Range("C6").Value = $enNew

# --- Column width adjustments ---
$ws.Columns.Item(2).ColumnWidth = 74.33333333333333
$ws.Columns.Item(3).ColumnWidth = 100.16666666666667

# --- Update selection / scroll position ---
$ws.Range("C6").Select()
